$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.244.38'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.427.95'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '414.19'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.12'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.89%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.46%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.725'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.140'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.70'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.78%  '
$ws.Range("E12").Value = '  +3.35%  '
$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000216'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +8.25%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.965.33'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.20%  '
$ws.Range("E15").Value = '  -0.44%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.51'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.440.06'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("E18").Value = '  +2.60%  '
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '62.281.85'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '467.33'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +5.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '90.93'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.41%  '
$ws.Range("E23").Value = '  +2.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.52'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +5.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.30'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '33.09'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.79'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.67'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '11.95'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.34%  '
$ws.Range("E31").Value = '  -3.22%  '
$ws.Range("E32").Value = '  -2.21%  '
$ws.Range("E33").Value = '  -1.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '40.76'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.91'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +10.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0487'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.06'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.92%  '
$ws.Range("E40").Value = '  +4.38%  '
$ws.Range("E41").Value = '  +0.31%  '
$ws.Range("E42").Value = '  -1.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.68'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +11.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '145.06'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.93%  '
$ws.Range("E45").Value = '  +5.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.31'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.46'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +19.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '16.41'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.29'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₃0525'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +30.54%  '
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '110.36'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +6.71%  '
